$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Java Full Stack AI Engineer (Senior Software Engineer)'
$ws.Cells.Item(2, 2).Value = 'LTIMindtree'
$ws.Cells.Item(2, 3).Value = 'Tampa, FL, US USA'
$ws.Cells.Item(2, 4).Value = 23.3
$ws.Cells.Item(2, 5).Value = 'AI Engineer, LangChain, RAG, LLaMA, Gemini, Copilot, Pinecone, Prompt Engineering, Kinesis, Docker'
$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = '2026-02-23'
$ws.Cells.Item(2, 7).Value = 'https://www.indeed.com/viewjob?jk=8915b5aa273f1b5b'

$ws.Cells.Item(3, 1).Value = 'AI Engineer'
$ws.Cells.Item(3, 2).Value = 'Mount Tech'
$ws.Cells.Item(3, 3).Value = 'La Jolla, CA, US USA'
$ws.Cells.Item(3, 4).Value = 21.1
$ws.Cells.Item(3, 5).Value = 'AI Engineer, LangChain, RAG, LLaMA, Copilot, Pinecone, Prompt Engineering, TensorFlow, PyTorch, FastAPI'
$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 6).Value = '2026-02-23'
$ws.Cells.Item(3, 7).Value = 'https://www.indeed.com/viewjob?jk=ba4cb3da56fe3cf3'

$ws.Cells.Item(4, 1).Value = 'Senior Data Engineer'
$ws.Cells.Item(4, 2).Value = 'Mariner Wealth Advisors'
$ws.Cells.Item(4, 3).Value = 'Overland Park, KS, US USA'
$ws.Cells.Item(4, 4).Value = 17.8
$ws.Cells.Item(4, 5).Value = 'RAG, Cortex, S3, Redshift, FastAPI, CI/CD, Terraform, Git, Snowflake, Databricks'
$ws.Cells.Item(4, 6).NumberFormat = "@"
$ws.Cells.Item(4, 6).Value = '2026-02-23'
$ws.Cells.Item(4, 7).Value = 'https://www.indeed.com/viewjob?jk=074ebeab3093055c'

$ws.Cells.Item(5, 1).Value = 'AI/Machine Learning Data Engineer'
$ws.Cells.Item(5, 2).Value = 'Acosta Group'
$ws.Cells.Item(5, 3).Value = 'Chesterfield, MO, US USA'
$ws.Cells.Item(5, 4).Value = 15.6
$ws.Cells.Item(5, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Copilot, Azure ML, Docker, Kubernetes, CI/CD, Databricks'
$ws.Cells.Item(5, 6).NumberFormat = "@"
$ws.Cells.Item(5, 6).Value = '2026-02-23'
$ws.Cells.Item(5, 7).Value = 'https://www.indeed.com/viewjob?jk=65d238209415bd35'

$ws.Cells.Item(6, 1).Value = 'AI/Machine Learning Data Engineer'
$ws.Cells.Item(6, 2).Value = 'Acosta Group'
$ws.Cells.Item(6, 3).Value = 'Jacksonville, FL, US USA'
$ws.Cells.Item(6, 4).Value = 15.6
$ws.Cells.Item(6, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Copilot, Azure ML, Docker, Kubernetes, CI/CD, Databricks'
$ws.Cells.Item(6, 6).NumberFormat = "@"
$ws.Cells.Item(6, 6).Value = '2026-02-23'
$ws.Cells.Item(6, 7).Value = 'https://www.indeed.com/viewjob?jk=08527261a275aa62'

$ws.Cells.Item(7, 1).Value = 'AI/Machine Learning Data Engineer'
$ws.Cells.Item(7, 2).Value = 'Acosta Group'
$ws.Cells.Item(7, 3).Value = 'Lewisville, TX, US USA'
$ws.Cells.Item(7, 4).Value = 15.6
$ws.Cells.Item(7, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Copilot, Azure ML, Docker, Kubernetes, CI/CD, Databricks'
$ws.Cells.Item(7, 6).NumberFormat = "@"
$ws.Cells.Item(7, 6).Value = '2026-02-23'
$ws.Cells.Item(7, 7).Value = 'https://www.indeed.com/viewjob?jk=dc66fc514f39b781'

$ws.Cells.Item(8, 1).Value = 'AI/Machine Learning Data Engineer'
$ws.Cells.Item(8, 2).Value = 'Acosta Group'
$ws.Cells.Item(8, 3).Value = 'Boise, ID, US USA'
$ws.Cells.Item(8, 4).Value = 15.6
$ws.Cells.Item(8, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Copilot, Azure ML, Docker, Kubernetes, CI/CD, Databricks'
$ws.Cells.Item(8, 6).NumberFormat = "@"
$ws.Cells.Item(8, 6).Value = '2026-02-23'
$ws.Cells.Item(8, 7).Value = 'https://www.indeed.com/viewjob?jk=003cc1fe2c6f4d82'

$ws.Cells.Item(9, 1).Value = 'AI/Machine Learning Data Engineer'
$ws.Cells.Item(9, 2).Value = 'Acosta Group'
$ws.Cells.Item(9, 3).Value = 'Chicago, IL, US USA'
$ws.Cells.Item(9, 4).Value = 15.6
$ws.Cells.Item(9, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Copilot, Azure ML, Docker, Kubernetes, CI/CD, Databricks'
$ws.Cells.Item(9, 6).NumberFormat = "@"
$ws.Cells.Item(9, 6).Value = '2026-02-23'
$ws.Cells.Item(9, 7).Value = 'https://www.indeed.com/viewjob?jk=ceeb3d3f63ddd02b'

$ws.Cells.Item(10, 1).Value = 'AI/Machine Learning Data Engineer'
$ws.Cells.Item(10, 2).Value = 'Acosta Group'
$ws.Cells.Item(10, 3).Value = 'Chesterfield, MO, US USA'
$ws.Cells.Item(10, 4).Value = 15.6
$ws.Cells.Item(10, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Copilot, Azure ML, Docker, Kubernetes, CI/CD, Databricks'
$ws.Cells.Item(10, 6).NumberFormat = "@"
$ws.Cells.Item(10, 6).Value = '2026-02-23'
$ws.Cells.Item(10, 7).Value = 'https://www.indeed.com/viewjob?jk=6394fdb64160187b'

$ws.Cells.Item(11, 1).Value = 'AI/Machine Learning Data Engineer'
$ws.Cells.Item(11, 2).Value = 'Acosta Group'
$ws.Cells.Item(11, 3).Value = 'Lewisville, TX, US USA'
$ws.Cells.Item(11, 4).Value = 15.6
$ws.Cells.Item(11, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Copilot, Azure ML, Docker, Kubernetes, CI/CD, Databricks'
$ws.Cells.Item(11, 6).NumberFormat = "@"
$ws.Cells.Item(11, 6).Value = '2026-02-23'
$ws.Cells.Item(11, 7).Value = 'https://www.indeed.com/viewjob?jk=4bf2ba0a6a2d2ef7'

$ws.Cells.Item(12, 1).Value = 'AI/Machine Learning Data Engineer'
$ws.Cells.Item(12, 2).Value = 'Acosta Group'
$ws.Cells.Item(12, 3).Value = 'Boise, ID, US USA'
$ws.Cells.Item(12, 4).Value = 15.6
$ws.Cells.Item(12, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Copilot, Azure ML, Docker, Kubernetes, CI/CD, Databricks'
$ws.Cells.Item(12, 6).NumberFormat = "@"
$ws.Cells.Item(12, 6).Value = '2026-02-23'
$ws.Cells.Item(12, 7).Value = 'https://www.indeed.com/viewjob?jk=6c579e60e29085a0'

$ws.Cells.Item(13, 1).Value = 'AI/Machine Learning Data Engineer'
$ws.Cells.Item(13, 2).Value = 'Acosta Group'
$ws.Cells.Item(13, 3).Value = 'Chicago, IL, US USA'
$ws.Cells.Item(13, 4).Value = 15.6
$ws.Cells.Item(13, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Copilot, Azure ML, Docker, Kubernetes, CI/CD, Databricks'
$ws.Cells.Item(13, 6).NumberFormat = "@"
$ws.Cells.Item(13, 6).Value = '2026-02-23'
$ws.Cells.Item(13, 7).Value = 'https://www.indeed.com/viewjob?jk=9566bb3277fc5b3e'

$ws.Cells.Item(14, 1).Value = 'AI/Machine Learning Data Engineer'
$ws.Cells.Item(14, 2).Value = 'Acosta Group'
$ws.Cells.Item(14, 3).Value = 'Rogers, AR, US USA'
$ws.Cells.Item(14, 4).Value = 15.6
$ws.Cells.Item(14, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Copilot, Azure ML, Docker, Kubernetes, CI/CD, Databricks'
$ws.Cells.Item(14, 6).NumberFormat = "@"
$ws.Cells.Item(14, 6).Value = '2026-02-23'
$ws.Cells.Item(14, 7).Value = 'https://www.indeed.com/viewjob?jk=8b988c6d60df206b'

$ws.Cells.Item(15, 1).Value = 'AI/Machine Learning Data Engineer'
$ws.Cells.Item(15, 2).Value = 'Acosta Group'
$ws.Cells.Item(15, 3).Value = 'Jacksonville, FL, US USA'
$ws.Cells.Item(15, 4).Value = 15.6
$ws.Cells.Item(15, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Copilot, Azure ML, Docker, Kubernetes, CI/CD, Databricks'
$ws.Cells.Item(15, 6).NumberFormat = "@"
$ws.Cells.Item(15, 6).Value = '2026-02-23'
$ws.Cells.Item(15, 7).Value = 'https://www.indeed.com/viewjob?jk=40ae88bd702dcb13'

$ws.Cells.Item(16, 1).Value = 'AI/Machine Learning Data Engineer'
$ws.Cells.Item(16, 2).Value = 'Acosta Group'
$ws.Cells.Item(16, 3).Value = 'Rogers, AR, US USA'
$ws.Cells.Item(16, 4).Value = 15.6
$ws.Cells.Item(16, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Copilot, Azure ML, Docker, Kubernetes, CI/CD, Databricks'
$ws.Cells.Item(16, 6).NumberFormat = "@"
$ws.Cells.Item(16, 6).Value = '2026-02-23'
$ws.Cells.Item(16, 7).Value = 'https://www.indeed.com/viewjob?jk=e490a596b2e349ab'

$ws.Cells.Item(17, 1).Value = 'Software Engineer II'
$ws.Cells.Item(17, 2).Value = 'Fanatics'
$ws.Cells.Item(17, 3).Value = 'US USA'
$ws.Cells.Item(17, 4).Value = 14.4
$ws.Cells.Item(17, 5).Value = 'Data Scientist, RAG, S3, Redshift, Kinesis, CI/CD, Git, Redshift, Kafka, MongoDB'
$ws.Cells.Item(17, 6).NumberFormat = "@"
$ws.Cells.Item(17, 6).Value = '2026-02-23'
$ws.Cells.Item(17, 7).Value = 'https://www.indeed.com/viewjob?jk=1be8fdec648472af'

$ws.Cells.Item(18, 1).Value = 'Senior Software Engineer - Generative AI'
$ws.Cells.Item(18, 2).Value = 'Acuity Insurance'
$ws.Cells.Item(18, 3).Value = 'Sheboygan, WI, US USA'
$ws.Cells.Item(18, 4).Value = 13.3
$ws.Cells.Item(18, 5).Value = 'Generative AI, LangChain, RAG, Prompt Engineering, CI/CD, PostgreSQL, Python, SQL, R, Java'
$ws.Cells.Item(18, 6).NumberFormat = "@"
$ws.Cells.Item(18, 6).Value = '2026-02-23'
$ws.Cells.Item(18, 7).Value = 'https://www.indeed.com/viewjob?jk=56a8c0e45b588390'

$ws.Cells.Item(19, 1).Value = 'Data Scientist'
$ws.Cells.Item(19, 2).Value = 'Pacific Community Ventures'
$ws.Cells.Item(19, 3).Value = 'Oakland, CA, US USA'
$ws.Cells.Item(19, 4).Value = 13.3
$ws.Cells.Item(19, 5).Value = 'Data Scientist, RAG, TensorFlow, PyTorch, CI/CD, Git, NoSQL, Tableau, Python, SQL'
$ws.Cells.Item(19, 6).NumberFormat = "@"
$ws.Cells.Item(19, 6).Value = '2026-02-23'
$ws.Cells.Item(19, 7).Value = 'https://www.indeed.com/viewjob?jk=0447d34a5a96e2c3'

$ws.Cells.Item(20, 1).Value = 'AI Engineer Mid-SR'
$ws.Cells.Item(20, 2).Value = 'Metova'
$ws.Cells.Item(20, 3).Value = 'PR, US USA'
$ws.Cells.Item(20, 4).Value = 13.3
$ws.Cells.Item(20, 5).Value = 'AI Engineer, LangChain, RAG, LLaMA, Pinecone, FastAPI, Docker, Kubernetes, CI/CD, Python'
$ws.Cells.Item(20, 6).NumberFormat = "@"
$ws.Cells.Item(20, 6).Value = '2026-02-23'
$ws.Cells.Item(20, 7).Value = 'https://www.indeed.com/viewjob?jk=accbd77f01ca5cfb'

$ws.Cells.Item(21, 1).Value = 'ML Ops Engineer II'
$ws.Cells.Item(21, 2).Value = 'Early Warning Services'
$ws.Cells.Item(21, 3).Value = 'San Francisco, CA, US USA'
$ws.Cells.Item(21, 4).Value = 13.3
$ws.Cells.Item(21, 5).Value = 'Data Scientist, RAG, MLflow, Docker, Kubernetes, CI/CD, Git, Hadoop, Python, R'
$ws.Cells.Item(21, 6).NumberFormat = "@"
$ws.Cells.Item(21, 6).Value = '2026-02-23'
$ws.Cells.Item(21, 7).Value = 'https://www.indeed.com/viewjob?jk=1bb70d7640bf48cc'

$ws.Cells.Item(22, 1).Value = 'Sr Machine Learning Engineer'
$ws.Cells.Item(22, 2).Value = 'The Walt Disney Company'
$ws.Cells.Item(22, 3).Value = 'Lake Buena Vista, FL, US USA'
$ws.Cells.Item(22, 4).Value = 13.3
$ws.Cells.Item(22, 5).Value = 'AI Engineer, Machine Learning Engineer, RAG, TensorFlow, PyTorch, Docker, Kubernetes, Git, Python, R'
$ws.Cells.Item(22, 6).NumberFormat = "@"
$ws.Cells.Item(22, 6).Value = '2026-02-23'
$ws.Cells.Item(22, 7).Value = 'https://www.indeed.com/viewjob?jk=fc65c03a732fa6c9'

$ws.Cells.Item(23, 1).Value = 'Enterprise Data Architect'
$ws.Cells.Item(23, 2).Value = 'IvoryCloud'
$ws.Cells.Item(23, 3).Value = 'Rockville, MD, US USA'
$ws.Cells.Item(23, 4).Value = 12.2
$ws.Cells.Item(23, 5).Value = 'RAG, Redshift, Data Lake, CI/CD, Snowflake, Databricks, Redshift, Python, SQL, R'
$ws.Cells.Item(23, 6).NumberFormat = "@"
$ws.Cells.Item(23, 6).Value = '2026-02-23'
$ws.Cells.Item(23, 7).Value = 'https://www.indeed.com/viewjob?jk=180dbcd7c9c036ee'

$ws.Cells.Item(24, 1).Value = 'Specialist, Data Engineer'
$ws.Cells.Item(24, 2).Value = 'Nationwide Mutual Insurance Company'
$ws.Cells.Item(24, 3).Value = 'Columbus, OH, US USA'
$ws.Cells.Item(24, 4).Value = 12.2
$ws.Cells.Item(24, 5).Value = 'Data Scientist, RAG, CI/CD, Jenkins, Git, Snowflake, Databricks, Python, SQL, R'
$ws.Cells.Item(24, 6).NumberFormat = "@"
$ws.Cells.Item(24, 6).Value = '2026-02-23'
$ws.Cells.Item(24, 7).Value = 'https://www.indeed.com/viewjob?jk=f15a9458813e7e74'

$ws.Cells.Item(25, 1).Value = 'Data Scientist'
$ws.Cells.Item(25, 2).Value = 'Indeed'
$ws.Cells.Item(25, 3).Value = 'Austin, TX, US USA'
$ws.Cells.Item(25, 4).Value = 12.2
$ws.Cells.Item(25, 5).Value = 'Data Scientist, RAG, Hadoop, Tableau, Power BI, Matplotlib, Seaborn, Python, SQL, R'
$ws.Cells.Item(25, 6).NumberFormat = "@"
$ws.Cells.Item(25, 6).Value = '2026-02-23'
$ws.Cells.Item(25, 7).Value = 'https://www.indeed.com/viewjob?jk=2ae3389e17a7f4d6'

$ws.Cells.Item(26, 1).Value = 'Applied AI Engineer'
$ws.Cells.Item(26, 2).Value = 'propio'
$ws.Cells.Item(26, 3).Value = 'Overland Park, KS, US USA'
$ws.Cells.Item(26, 4).Value = 12.2
$ws.Cells.Item(26, 5).Value = 'AI Engineer, LangChain, Hugging Face, FAISS, Pinecone, Prompt Engineering, FastAPI, Python, R, Scala'
$ws.Cells.Item(26, 6).NumberFormat = "@"
$ws.Cells.Item(26, 6).Value = '2026-02-23'
$ws.Cells.Item(26, 7).Value = 'https://www.indeed.com/viewjob?jk=223e2d3d1d722fdd'

$ws.Cells.Item(27, 1).Value = 'Solutions Architect'
$ws.Cells.Item(27, 2).Value = 'Interworks'
$ws.Cells.Item(27, 3).Value = 'Oklahoma City, OK, US USA'
$ws.Cells.Item(27, 4).Value = 11.1
$ws.Cells.Item(27, 5).Value = 'Glue, Redshift, BigQuery, Snowflake, Databricks, BigQuery, Redshift, Python, R, Scala'
$ws.Cells.Item(27, 6).NumberFormat = "@"
$ws.Cells.Item(27, 6).Value = '2026-02-23'
$ws.Cells.Item(27, 7).Value = 'https://www.indeed.com/viewjob?jk=ae3df7fdbe227203'

$ws.Cells.Item(28, 1).Value = 'Solutions Architect'
$ws.Cells.Item(28, 2).Value = 'Interworks'
$ws.Cells.Item(28, 3).Value = 'Stillwater, OK, US USA'
$ws.Cells.Item(28, 4).Value = 11.1
$ws.Cells.Item(28, 5).Value = 'Glue, Redshift, BigQuery, Snowflake, Databricks, BigQuery, Redshift, Python, R, Scala'
$ws.Cells.Item(28, 6).NumberFormat = "@"
$ws.Cells.Item(28, 6).Value = '2026-02-23'
$ws.Cells.Item(28, 7).Value = 'https://www.indeed.com/viewjob?jk=b2227ed8fd327afa'

$ws.Cells.Item(29, 1).Value = 'Solutions Architect'
$ws.Cells.Item(29, 2).Value = 'Interworks'
$ws.Cells.Item(29, 3).Value = 'Raleigh, NC, US USA'
$ws.Cells.Item(29, 4).Value = 11.1
$ws.Cells.Item(29, 5).Value = 'Glue, Redshift, BigQuery, Snowflake, Databricks, BigQuery, Redshift, Python, R, Scala'
$ws.Cells.Item(29, 6).NumberFormat = "@"
$ws.Cells.Item(29, 6).Value = '2026-02-23'
$ws.Cells.Item(29, 7).Value = 'https://www.indeed.com/viewjob?jk=a40f1918af374ac1'

$ws.Cells.Item(30, 1).Value = 'Solutions Architect'
$ws.Cells.Item(30, 2).Value = 'Interworks'
$ws.Cells.Item(30, 3).Value = 'Tulsa, OK, US USA'
$ws.Cells.Item(30, 4).Value = 11.1
$ws.Cells.Item(30, 5).Value = 'Glue, Redshift, BigQuery, Snowflake, Databricks, BigQuery, Redshift, Python, R, Scala'
$ws.Cells.Item(30, 6).NumberFormat = "@"
$ws.Cells.Item(30, 6).Value = '2026-02-23'
$ws.Cells.Item(30, 7).Value = 'https://www.indeed.com/viewjob?jk=b65ec0cc7ed25eda'

$ws.Cells.Item(31, 1).Value = 'Risk Adjustment Sr. Data Analyst - Remote'
$ws.Cells.Item(31, 2).Value = 'Datavant'
$ws.Cells.Item(31, 3).Value = 'Houston, TX, US USA'
$ws.Cells.Item(31, 4).Value = 10
$ws.Cells.Item(31, 5).Value = 'RAG, Snowflake, Databricks, Tableau, Power BI, Python, SQL, R, Scala'
$ws.Cells.Item(31, 6).NumberFormat = "@"
$ws.Cells.Item(31, 6).Value = '2026-02-23'
$ws.Cells.Item(31, 7).Value = 'https://www.indeed.com/viewjob?jk=6d255eff083655ae'

$ws.Cells.Item(32, 1).Value = 'Data Scientist/AI Trainer'
$ws.Cells.Item(32, 2).Value = 'Five9'
$ws.Cells.Item(32, 3).Value = 'Remote, US USA'
$ws.Cells.Item(32, 4).Value = 10
$ws.Cells.Item(32, 5).Value = 'Data Scientist, RAG, Prompt Engineering, TensorFlow, PyTorch, Python, SQL, R, Scala'
$ws.Cells.Item(32, 6).NumberFormat = "@"
$ws.Cells.Item(32, 6).Value = '2026-02-23'
$ws.Cells.Item(32, 7).Value = 'https://www.indeed.com/viewjob?jk=0d85ec88c21eba2c'

$ws.Cells.Item(33, 1).Value = 'Senior QA Engineer'
$ws.Cells.Item(33, 2).Value = 'Care.com'
$ws.Cells.Item(33, 3).Value = 'Salt Lake City, UT, US USA'
$ws.Cells.Item(33, 4).Value = 10
$ws.Cells.Item(33, 5).Value = 'RAG, CI/CD, Jenkins, GitHub Actions, Git, Python, R, Java, Scala'
$ws.Cells.Item(33, 6).NumberFormat = "@"
$ws.Cells.Item(33, 6).Value = '2026-02-23'
$ws.Cells.Item(33, 7).Value = 'https://www.indeed.com/viewjob?jk=22b669fc82caea43'

$ws.Cells.Item(34, 1).Value = 'Senior QA Engineer'
$ws.Cells.Item(34, 2).Value = 'Care.com'
$ws.Cells.Item(34, 3).Value = 'Dallas, TX, US USA'
$ws.Cells.Item(34, 4).Value = 10
$ws.Cells.Item(34, 5).Value = 'RAG, CI/CD, Jenkins, GitHub Actions, Git, Python, R, Java, Scala'
$ws.Cells.Item(34, 6).NumberFormat = "@"
$ws.Cells.Item(34, 6).Value = '2026-02-23'
$ws.Cells.Item(34, 7).Value = 'https://www.indeed.com/viewjob?jk=4acc2278c1957139'

$ws.Cells.Item(35, 1).Value = 'Senior QA Engineer'
$ws.Cells.Item(35, 2).Value = 'Care.com'
$ws.Cells.Item(35, 3).Value = 'Austin, TX, US USA'
$ws.Cells.Item(35, 4).Value = 10
$ws.Cells.Item(35, 5).Value = 'RAG, CI/CD, Jenkins, GitHub Actions, Git, Python, R, Java, Scala'
$ws.Cells.Item(35, 6).NumberFormat = "@"
$ws.Cells.Item(35, 6).Value = '2026-02-23'
$ws.Cells.Item(35, 7).Value = 'https://www.indeed.com/viewjob?jk=9870ae11a2356e86'

$ws.Cells.Item(36, 1).Value = 'Architect, Service & Operational Data'
$ws.Cells.Item(36, 2).Value = 'Thomson Reuters'
$ws.Cells.Item(36, 3).Value = 'Eagan, MN, US USA'
$ws.Cells.Item(36, 4).Value = 10
$ws.Cells.Item(36, 5).Value = 'RAG, Kinesis, CI/CD, Git, Kafka, Python, SQL, R, Scala'
$ws.Cells.Item(36, 6).NumberFormat = "@"
$ws.Cells.Item(36, 6).Value = '2026-02-23'
$ws.Cells.Item(36, 7).Value = 'https://www.indeed.com/viewjob?jk=c82df1ac5ad3ee0a'

$ws.Cells.Item(37, 1).Value = 'Data Scientist'
$ws.Cells.Item(37, 2).Value = 'Ipsos'
$ws.Cells.Item(37, 3).Value = 'Culver City, CA, US USA'
$ws.Cells.Item(37, 4).Value = 10
$ws.Cells.Item(37, 5).Value = 'Data Scientist, RAG, Git, Python, SQL, R, Java, Scala, Bayesian'
$ws.Cells.Item(37, 6).NumberFormat = "@"
$ws.Cells.Item(37, 6).Value = '2026-02-23'
$ws.Cells.Item(37, 7).Value = 'https://www.indeed.com/viewjob?jk=7b2acc4634dfe573'

$ws.Cells.Item(38, 1).Value = 'Senior Big Data Engineer'
$ws.Cells.Item(38, 2).Value = 'Highmark Health'
$ws.Cells.Item(38, 3).Value = 'PA, US USA'
$ws.Cells.Item(38, 4).Value = 10
$ws.Cells.Item(38, 5).Value = 'AI Engineer, Data Scientist, Git, Kafka, Tableau, Python, SQL, R, Scala'
$ws.Cells.Item(38, 6).NumberFormat = "@"
$ws.Cells.Item(38, 6).Value = '2026-02-23'
$ws.Cells.Item(38, 7).Value = 'https://www.indeed.com/viewjob?jk=e566377201b9121f'

$ws.Cells.Item(39, 1).Value = 'Security Software Developer'
$ws.Cells.Item(39, 2).Value = 'Vidoori'
$ws.Cells.Item(39, 3).Value = 'Hyattsville, MD, US USA'
$ws.Cells.Item(39, 4).Value = 10
$ws.Cells.Item(39, 5).Value = 'RAG, Docker, Kubernetes, CI/CD, Git, Python, R, Java, Scala'
$ws.Cells.Item(39, 6).NumberFormat = "@"
$ws.Cells.Item(39, 6).Value = '2026-02-23'
$ws.Cells.Item(39, 7).Value = 'https://www.indeed.com/viewjob?jk=073709b581c90d36'

